# Update row 2 (the single data record) of Sheet1 to reflect the new
# patient / requisition / insurance data, and re-shuffle which columns
# are populated (MiddleName + PrimaryInsurance_SubDOB/ContractNumber are
# now populated while the SecondaryInsurance_* columns are cleared).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $rng = $ws.Range($CellRef)
    # Force a Text number format before assigning so that values which look
    # like dates/numbers (e.g. "2024-02-14", "33514") are kept as literal
    # strings instead of being parsed into serial date/number values.
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    # Drop the explicit Text format again so the cell keeps using the
    # workbook's default style (style index 0), matching the rest of the
    # sheet; the stored value remains the literal text already assigned.
    $rng.ClearFormats()
}

Set-TextValue "A2" "Third-Party Bill"
Set-TextValue "B2" "2024-02-14"
Set-TextValue "C2" "33514"
Set-TextValue "D2" "29469205"
Set-TextValue "E2" "SEWMA001"
Set-TextValue "F2" "Sagis DX"
Set-TextValue "G2" "SEWARD"
Set-TextValue "H2" "MARY"
Set-TextValue "I2" "B"
Set-TextValue "J2" "1954-10-20"
Set-TextValue "K2" "Female"
Set-TextValue "L2" "826 AUTUMN PL"
Set-TextValue "M2" "LA"
Set-TextValue "N2" "MANDEVILLE"
Set-TextValue "O2" "704716772"
Set-TextValue "P2" "5044603163"
Set-TextValue "Q2" "table"
Set-TextValue "R2" "Stratton Beatrous Grisoli M.D.,"
Set-TextValue "S2" "Baldone Reina Dermatology, APMC"
Set-TextValue "T2" "MARY SEWARD"
Set-TextValue "U2" "Self"
Set-TextValue "V2" "1954-10-20"
Set-TextValue "W2" "Humana Health Plan"
Set-TextValue "X2" "X2284001"
Set-TextValue "Y2" "H64587383"

# The SecondaryInsurance_* fields (SubscriberName, Relationship,
# CompanyName, MemberID) are no longer populated for this record.
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AG2").ClearContents()
$ws.Range("AI2").ClearContents()
